$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price/volume updates ---
$ws.Range('D2').Value = '67.496.47'
$ws.Range('E2').Value = '  -1.42%  '
$ws.Range('D3').Value = '3.322.85'
$ws.Range('E3').Value = '  -1.39%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '582.41'
$ws.Range('E5').Value = '  -1.89%  '
$ws.Range('D6').Value = '175.69'
$ws.Range('E6').Value = '  -5.31%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('E8').Value = '  -1.34%  '
$ws.Range('D9').Value = '3.320.25'
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('E10').Value = '  -2.87%  '
$ws.Range('D12').Value = '45.38'
$ws.Range('E12').Value = '  -3.62%  '
$ws.Range('E13').Value = '  -3.65%  '
$ws.Range('D14').Value = '658.99'
$ws.Range('E14').Value = '  +2.91%  '
$ws.Range('D15').Value = '3.866.03'
$ws.Range('E15').Value = '  -1.06%  '
$ws.Range('D16').Value = '8.38'
$ws.Range('E16').Value = '  -1.70%  '
$ws.Range('D17').Value = '67.671.26'
$ws.Range('E17').Value = '  -1.51%  '
$ws.Range('D19').Value = '3.325.74'
$ws.Range('E19').Value = '  -1.38%  '
$ws.Range('D21').Value = '10.91'
$ws.Range('E21').Value = '  -1.17%  '
$ws.Range('E22').Value = '  -2.34%  '
$ws.Range('D23').Value = '5.43'
$ws.Range('E23').Value = '  +7.06%  '
$ws.Range('E24').Value = '  -5.04%  '
$ws.Range('D25').Value = '99.27'
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('D26').Value = '3.84'
$ws.Range('E26').Value = '  -6.03%  '
$ws.Range('E27').Value = '  -6.14%  '
$ws.Range('E28').Value = '  -5.35%  '
$ws.Range('D29').Value = '33.54'
$ws.Range('E29').Value = '  +1.78%  '
$ws.Range('D30').Value = '7.40'
$ws.Range('E30').Value = '  +9.04%  '
$ws.Range('D31').Value = '8.41'
$ws.Range('E31').Value = '  -2.82%  '
$ws.Range('D32').Value = '591.39'
$ws.Range('E32').Value = '  -3.07%  '
$ws.Range('D33').Value = '10.96'
$ws.Range('E33').Value = '  -0.88%  '
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').Value = '3.695.21'
$ws.Range('E36').Value = '  -7.27%  '
$ws.Range('D37').Value = '56.52'
$ws.Range('E37').Value = '  +0.77%  '
$ws.Range('D38').Value = '3.34'
$ws.Range('E38').Value = '  -9.17%  '
$ws.Range('E41').Value = '  -4.93%  '
$ws.Range('E42').Value = '  -5.89%  '
$ws.Range('D43').Value = '0.332'
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('E44').Value = '  -5.34%  '
$ws.Range('D45').Value = '3.26'
$ws.Range('E45').Value = '  -4.96%  '
$ws.Range('E46').Value = '  -3.42%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('E50').Value = '  -0.21%  '
$ws.Range('D51').Value = '127.07'
$ws.Range('E51').Value = '  -3.23%  '

# --- Row swaps (coin re-ranking) ---
$ws.Range('B39').Value = 'InjectiveProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D39').Value = '33.98'
$ws.Range('E39').Value = '  +1.49%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').Value = '0.130'
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '0.128'
$ws.Range('E47').Value = '  -1.21%  '
$ws.Range('B48').Value = 'ThetaToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D48').Value = '2.58'
$ws.Range('E48').Value = '  +0.12%  '

Write-Output "Updated cryptos list"